$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header strings: _old -> _FV2210, _new -> _FV2304
$newHeaders = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($c = 1; $c -le 21; $c++) {
    $ws.Cells.Item(1, $c).Value = $newHeaders[$c - 1]
}

# 2. Add table over A1:U58
$range = $ws.Range("A1:U58")
$xlSrcRange = [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange
$xlYes = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
$listObject = $ws.ListObjects.Add($xlSrcRange, $range, $null, $xlYes)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# 3. Freeze top row (pane split)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
